$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# In the source workbook, the two shared strings used for the "Requisitos"
# rows (24 and 25) had their order swapped. Row/column mapping in the sheet
# itself did not change, so the net visible effect is that the text shown
# in row 24 and row 25 (columns B and C) trade places:
#   row 24 (was "LOM3202 -  Circuitos Elétricos  (Indicação de Conjunto)")
#          becomes "LOB1053 -  Física III  (Requisito)"
#   row 25 (was "LOB1053 -  Física III  (Requisito)")
#          becomes "LOM3202 -  Circuitos Elétricos  (Indicação de Conjunto)"

$row24Text = $ws.Range("B24").Value()
$row25Text = $ws.Range("B25").Value()

$ws.Range("B24").Value = $row25Text
$ws.Range("C24").Value = $row25Text
$ws.Range("B25").Value = $row24Text
$ws.Range("C25").Value = $row24Text
